# Natmi following Dr Hou advice
# Adds the "M2" sending-cluster block (rows 4,8,12 D-values + new rows 14-17)
# and refreshes all Sema6d-Kdr LR-pair statistics for Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Sema6d"
$ws.Cells.Item(2,3).Value = "Kdr"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = [double]"3"
$ws.Cells.Item(2,6).Value = [double]"1"
$ws.Cells.Item(2,7).Value = [double]"24.44575933333333"
$ws.Cells.Item(2,8).Value = [double]"73.337278"
$ws.Cells.Item(2,9).Value = [double]"0.4034052273345712"
$ws.Cells.Item(2,10).Value = [double]"0.4034052273345712"
$ws.Cells.Item(2,11).Value = [double]"3"
$ws.Cells.Item(2,12).Value = [double]"1"
$ws.Cells.Item(2,13).Value = [double]"127.6999736666667"
$ws.Cells.Item(2,14).Value = [double]"383.099921"
$ws.Cells.Item(2,15).Value = [double]"0.9554352891750322"
$ws.Cells.Item(2,16).Value = [double]"0.9554352891750322"
$ws.Cells.Item(2,17).Value = [double]"3121.722823128337"
$ws.Cells.Item(2,18).Value = [double]"28095.50540815504"
$ws.Cells.Item(2,19).Value = [double]"0.3854275900331257"
$ws.Cells.Item(2,20).Value = [double]"0.3854275900331256"

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Sema6d"
$ws.Cells.Item(3,3).Value = "Kdr"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = [double]"3"
$ws.Cells.Item(3,6).Value = [double]"1"
$ws.Cells.Item(3,7).Value = [double]"24.44575933333333"
$ws.Cells.Item(3,8).Value = [double]"73.337278"
$ws.Cells.Item(3,9).Value = [double]"0.4034052273345712"
$ws.Cells.Item(3,10).Value = [double]"0.4034052273345712"
$ws.Cells.Item(3,11).Value = [double]"3"
$ws.Cells.Item(3,12).Value = [double]"1"
$ws.Cells.Item(3,13).Value = [double]"0.4321196666666667"
$ws.Cells.Item(3,14).Value = [double]"1.296359"
$ws.Cells.Item(3,15).Value = [double]"0.003233065495828321"
$ws.Cells.Item(3,16).Value = [double]"0.003233065495828321"
$ws.Cells.Item(3,17).Value = [double]"10.56349337453356"
$ws.Cells.Item(3,18).Value = [double]"95.07144037080199"
$ws.Cells.Item(3,19).Value = [double]"0.001304235521332182"
$ws.Cells.Item(3,20).Value = [double]"0.001304235521332182"

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Sema6d"
$ws.Cells.Item(4,3).Value = "Kdr"
$ws.Cells.Item(4,4).Value = "M2"
$ws.Cells.Item(4,5).Value = [double]"3"
$ws.Cells.Item(4,6).Value = [double]"1"
$ws.Cells.Item(4,7).Value = [double]"24.44575933333333"
$ws.Cells.Item(4,8).Value = [double]"73.337278"
$ws.Cells.Item(4,9).Value = [double]"0.4034052273345712"
$ws.Cells.Item(4,10).Value = [double]"0.4034052273345712"
$ws.Cells.Item(4,11).Value = [double]"3"
$ws.Cells.Item(4,12).Value = [double]"1"
$ws.Cells.Item(4,13).Value = [double]"4.77305"
$ws.Cells.Item(4,14).Value = [double]"14.31915"
$ws.Cells.Item(4,15).Value = [double]"0.03571136528892854"
$ws.Cells.Item(4,16).Value = [double]"0.03571136528892854"
$ws.Cells.Item(4,17).Value = [double]"116.6808315859667"
$ws.Cells.Item(4,18).Value = [double]"1050.1274842737"
$ws.Cells.Item(4,19).Value = [double]"0.01440615143280813"
$ws.Cells.Item(4,20).Value = [double]"0.01440615143280813"

# Row 5
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Sema6d"
$ws.Cells.Item(5,3).Value = "Kdr"
$ws.Cells.Item(5,4).Value = "sCs"
$ws.Cells.Item(5,5).Value = [double]"3"
$ws.Cells.Item(5,6).Value = [double]"1"
$ws.Cells.Item(5,7).Value = [double]"24.44575933333333"
$ws.Cells.Item(5,8).Value = [double]"73.337278"
$ws.Cells.Item(5,9).Value = [double]"0.4034052273345712"
$ws.Cells.Item(5,10).Value = [double]"0.4034052273345712"
$ws.Cells.Item(5,11).Value = [double]"3"
$ws.Cells.Item(5,12).Value = [double]"1"
$ws.Cells.Item(5,13).Value = [double]"0.751186"
$ws.Cells.Item(5,14).Value = [double]"2.253558"
$ws.Cells.Item(5,15).Value = [double]"0.00562028004021099"
$ws.Cells.Item(5,16).Value = [double]"0.00562028004021099"
$ws.Cells.Item(5,17).Value = [double]"18.36331217056933"
$ws.Cells.Item(5,18).Value = [double]"165.269809535124"
$ws.Cells.Item(5,19).Value = [double]"0.002267250347305268"
$ws.Cells.Item(5,20).Value = [double]"0.002267250347305267"

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Sema6d"
$ws.Cells.Item(6,3).Value = "Kdr"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = [double]"3"
$ws.Cells.Item(6,6).Value = [double]"1"
$ws.Cells.Item(6,7).Value = [double]"14.28901333333333"
$ws.Cells.Item(6,8).Value = [double]"42.86704"
$ws.Cells.Item(6,9).Value = [double]"0.2357980618855278"
$ws.Cells.Item(6,10).Value = [double]"0.2357980618855278"
$ws.Cells.Item(6,11).Value = [double]"3"
$ws.Cells.Item(6,12).Value = [double]"1"
$ws.Cells.Item(6,13).Value = [double]"127.6999736666667"
$ws.Cells.Item(6,14).Value = [double]"383.099921"
$ws.Cells.Item(6,15).Value = [double]"0.9554352891750322"
$ws.Cells.Item(6,16).Value = [double]"0.9554352891750322"
$ws.Cells.Item(6,17).Value = [double]"1824.706626389316"
$ws.Cells.Item(6,18).Value = [double]"16422.35963750384"
$ws.Cells.Item(6,19).Value = [double]"0.2252897894445114"
$ws.Cells.Item(6,20).Value = [double]"0.2252897894445114"

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Sema6d"
$ws.Cells.Item(7,3).Value = "Kdr"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = [double]"3"
$ws.Cells.Item(7,6).Value = [double]"1"
$ws.Cells.Item(7,7).Value = [double]"14.28901333333333"
$ws.Cells.Item(7,8).Value = [double]"42.86704"
$ws.Cells.Item(7,9).Value = [double]"0.2357980618855278"
$ws.Cells.Item(7,10).Value = [double]"0.2357980618855278"
$ws.Cells.Item(7,11).Value = [double]"3"
$ws.Cells.Item(7,12).Value = [double]"1"
$ws.Cells.Item(7,13).Value = [double]"0.4321196666666667"
$ws.Cells.Item(7,14).Value = [double]"1.296359"
$ws.Cells.Item(7,15).Value = [double]"0.003233065495828321"
$ws.Cells.Item(7,16).Value = [double]"0.003233065495828321"
$ws.Cells.Item(7,17).Value = [double]"6.174563678595557"
$ws.Cells.Item(7,18).Value = [double]"55.57107310736001"
$ws.Cells.Item(7,19).Value = [double]"0.0007623505778652912"
$ws.Cells.Item(7,20).Value = [double]"0.0007623505778652911"

# Row 8
$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Sema6d"
$ws.Cells.Item(8,3).Value = "Kdr"
$ws.Cells.Item(8,4).Value = "M2"
$ws.Cells.Item(8,5).Value = [double]"3"
$ws.Cells.Item(8,6).Value = [double]"1"
$ws.Cells.Item(8,7).Value = [double]"14.28901333333333"
$ws.Cells.Item(8,8).Value = [double]"42.86704"
$ws.Cells.Item(8,9).Value = [double]"0.2357980618855278"
$ws.Cells.Item(8,10).Value = [double]"0.2357980618855278"
$ws.Cells.Item(8,11).Value = [double]"3"
$ws.Cells.Item(8,12).Value = [double]"1"
$ws.Cells.Item(8,13).Value = [double]"4.77305"
$ws.Cells.Item(8,14).Value = [double]"14.31915"
$ws.Cells.Item(8,15).Value = [double]"0.03571136528892854"
$ws.Cells.Item(8,16).Value = [double]"0.03571136528892854"
$ws.Cells.Item(8,17).Value = [double]"68.20217509066669"
$ws.Cells.Item(8,18).Value = [double]"613.8195758160001"
$ws.Cells.Item(8,19).Value = [double]"0.008420670722415462"
$ws.Cells.Item(8,20).Value = [double]"0.008420670722415461"

# Row 9
$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Sema6d"
$ws.Cells.Item(9,3).Value = "Kdr"
$ws.Cells.Item(9,4).Value = "sCs"
$ws.Cells.Item(9,5).Value = [double]"3"
$ws.Cells.Item(9,6).Value = [double]"1"
$ws.Cells.Item(9,7).Value = [double]"14.28901333333333"
$ws.Cells.Item(9,8).Value = [double]"42.86704"
$ws.Cells.Item(9,9).Value = [double]"0.2357980618855278"
$ws.Cells.Item(9,10).Value = [double]"0.2357980618855278"
$ws.Cells.Item(9,11).Value = [double]"3"
$ws.Cells.Item(9,12).Value = [double]"1"
$ws.Cells.Item(9,13).Value = [double]"0.751186"
$ws.Cells.Item(9,14).Value = [double]"2.253558"
$ws.Cells.Item(9,15).Value = [double]"0.00562028004021099"
$ws.Cells.Item(9,16).Value = [double]"0.00562028004021099"
$ws.Cells.Item(9,17).Value = [double]"10.73370676981333"
$ws.Cells.Item(9,18).Value = [double]"96.60336092832"
$ws.Cells.Item(9,19).Value = [double]"0.001325251140735668"
$ws.Cells.Item(9,20).Value = [double]"0.001325251140735668"

# Row 10
$ws.Cells.Item(10,1).Value = "M2"
$ws.Cells.Item(10,2).Value = "Sema6d"
$ws.Cells.Item(10,3).Value = "Kdr"
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(10,5).Value = [double]"2"
$ws.Cells.Item(10,6).Value = [double]"0.6666666666666666"
$ws.Cells.Item(10,7).Value = [double]"0.480005"
$ws.Cells.Item(10,8).Value = [double]"1.440015"
$ws.Cells.Item(10,9).Value = [double]"0.007921068169999337"
$ws.Cells.Item(10,10).Value = [double]"0.007921068169999336"
$ws.Cells.Item(10,11).Value = [double]"3"
$ws.Cells.Item(10,12).Value = [double]"1"
$ws.Cells.Item(10,13).Value = [double]"127.6999736666667"
$ws.Cells.Item(10,14).Value = [double]"383.099921"
$ws.Cells.Item(10,15).Value = [double]"0.9554352891750322"
$ws.Cells.Item(10,16).Value = [double]"0.9554352891750322"
$ws.Cells.Item(10,17).Value = [double]"61.29662585986834"
$ws.Cells.Item(10,18).Value = [double]"551.669632738815"
$ws.Cells.Item(10,19).Value = [double]"0.00756806805757846"
$ws.Cells.Item(10,20).Value = [double]"0.007568068057578458"

# Row 11
$ws.Cells.Item(11,1).Value = "M2"
$ws.Cells.Item(11,2).Value = "Sema6d"
$ws.Cells.Item(11,3).Value = "Kdr"
$ws.Cells.Item(11,4).Value = "FAPs"
$ws.Cells.Item(11,5).Value = [double]"2"
$ws.Cells.Item(11,6).Value = [double]"0.6666666666666666"
$ws.Cells.Item(11,7).Value = [double]"0.480005"
$ws.Cells.Item(11,8).Value = [double]"1.440015"
$ws.Cells.Item(11,9).Value = [double]"0.007921068169999337"
$ws.Cells.Item(11,10).Value = [double]"0.007921068169999336"
$ws.Cells.Item(11,11).Value = [double]"3"
$ws.Cells.Item(11,12).Value = [double]"1"
$ws.Cells.Item(11,13).Value = [double]"0.4321196666666667"
$ws.Cells.Item(11,14).Value = [double]"1.296359"
$ws.Cells.Item(11,15).Value = [double]"0.003233065495828321"
$ws.Cells.Item(11,16).Value = [double]"0.003233065495828321"
$ws.Cells.Item(11,17).Value = [double]"0.2074196005983333"
$ws.Cells.Item(11,18).Value = [double]"1.866776405385"
$ws.Cells.Item(11,19).Value = [double]"2.560933219052884E-05"
$ws.Cells.Item(11,20).Value = [double]"2.560933219052883E-05"

# Row 12
$ws.Cells.Item(12,1).Value = "M2"
$ws.Cells.Item(12,2).Value = "Sema6d"
$ws.Cells.Item(12,3).Value = "Kdr"
$ws.Cells.Item(12,4).Value = "M2"
$ws.Cells.Item(12,5).Value = [double]"2"
$ws.Cells.Item(12,6).Value = [double]"0.6666666666666666"
$ws.Cells.Item(12,7).Value = [double]"0.480005"
$ws.Cells.Item(12,8).Value = [double]"1.440015"
$ws.Cells.Item(12,9).Value = [double]"0.007921068169999337"
$ws.Cells.Item(12,10).Value = [double]"0.007921068169999336"
$ws.Cells.Item(12,11).Value = [double]"3"
$ws.Cells.Item(12,12).Value = [double]"1"
$ws.Cells.Item(12,13).Value = [double]"4.77305"
$ws.Cells.Item(12,14).Value = [double]"14.31915"
$ws.Cells.Item(12,15).Value = [double]"0.03571136528892854"
$ws.Cells.Item(12,16).Value = [double]"0.03571136528892854"
$ws.Cells.Item(12,17).Value = [double]"2.29108786525"
$ws.Cells.Item(12,18).Value = [double]"20.61979078725"
$ws.Cells.Item(12,19).Value = [double]"0.000282872158897351"
$ws.Cells.Item(12,20).Value = [double]"0.0002828721588973509"

# Row 13
$ws.Cells.Item(13,1).Value = "M2"
$ws.Cells.Item(13,2).Value = "Sema6d"
$ws.Cells.Item(13,3).Value = "Kdr"
$ws.Cells.Item(13,4).Value = "sCs"
$ws.Cells.Item(13,5).Value = [double]"2"
$ws.Cells.Item(13,6).Value = [double]"0.6666666666666666"
$ws.Cells.Item(13,7).Value = [double]"0.480005"
$ws.Cells.Item(13,8).Value = [double]"1.440015"
$ws.Cells.Item(13,9).Value = [double]"0.007921068169999337"
$ws.Cells.Item(13,10).Value = [double]"0.007921068169999336"
$ws.Cells.Item(13,11).Value = [double]"3"
$ws.Cells.Item(13,12).Value = [double]"1"
$ws.Cells.Item(13,13).Value = [double]"0.751186"
$ws.Cells.Item(13,14).Value = [double]"2.253558"
$ws.Cells.Item(13,15).Value = [double]"0.00562028004021099"
$ws.Cells.Item(13,16).Value = [double]"0.00562028004021099"
$ws.Cells.Item(13,17).Value = [double]"0.36057303593"
$ws.Cells.Item(13,18).Value = [double]"3.24515732337"
$ws.Cells.Item(13,19).Value = [double]"4.451862133299787E-05"
$ws.Cells.Item(13,20).Value = [double]"4.451862133299786E-05"

# Row 14
$ws.Cells.Item(14,1).Value = "sCs"
$ws.Cells.Item(14,2).Value = "Sema6d"
$ws.Cells.Item(14,3).Value = "Kdr"
$ws.Cells.Item(14,4).Value = "ECs"
$ws.Cells.Item(14,5).Value = [double]"3"
$ws.Cells.Item(14,6).Value = [double]"1"
$ws.Cells.Item(14,7).Value = [double]"21.38374133333333"
$ws.Cells.Item(14,8).Value = [double]"64.151224"
$ws.Cells.Item(14,9).Value = [double]"0.3528756426099016"
$ws.Cells.Item(14,10).Value = [double]"0.3528756426099016"
$ws.Cells.Item(14,11).Value = [double]"3"
$ws.Cells.Item(14,12).Value = [double]"1"
$ws.Cells.Item(14,13).Value = [double]"127.6999736666667"
$ws.Cells.Item(14,14).Value = [double]"383.099921"
$ws.Cells.Item(14,15).Value = [double]"0.9554352891750322"
$ws.Cells.Item(14,16).Value = [double]"0.9554352891750322"
$ws.Cells.Item(14,17).Value = [double]"2730.703205161478"
$ws.Cells.Item(14,18).Value = [double]"24576.3288464533"
$ws.Cells.Item(14,19).Value = [double]"0.3371498416398166"
$ws.Cells.Item(14,20).Value = [double]"0.3371498416398166"

# Row 15
$ws.Cells.Item(15,1).Value = "sCs"
$ws.Cells.Item(15,2).Value = "Sema6d"
$ws.Cells.Item(15,3).Value = "Kdr"
$ws.Cells.Item(15,4).Value = "FAPs"
$ws.Cells.Item(15,5).Value = [double]"3"
$ws.Cells.Item(15,6).Value = [double]"1"
$ws.Cells.Item(15,7).Value = [double]"21.38374133333333"
$ws.Cells.Item(15,8).Value = [double]"64.151224"
$ws.Cells.Item(15,9).Value = [double]"0.3528756426099016"
$ws.Cells.Item(15,10).Value = [double]"0.3528756426099016"
$ws.Cells.Item(15,11).Value = [double]"3"
$ws.Cells.Item(15,12).Value = [double]"1"
$ws.Cells.Item(15,13).Value = [double]"0.4321196666666667"
$ws.Cells.Item(15,14).Value = [double]"1.296359"
$ws.Cells.Item(15,15).Value = [double]"0.003233065495828321"
$ws.Cells.Item(15,16).Value = [double]"0.003233065495828321"
$ws.Cells.Item(15,17).Value = [double]"9.240335177046223"
$ws.Cells.Item(15,18).Value = [double]"83.163016593416"
$ws.Cells.Item(15,19).Value = [double]"0.001140870064440319"
$ws.Cells.Item(15,20).Value = [double]"0.001140870064440319"

# Row 16
$ws.Cells.Item(16,1).Value = "sCs"
$ws.Cells.Item(16,2).Value = "Sema6d"
$ws.Cells.Item(16,3).Value = "Kdr"
$ws.Cells.Item(16,4).Value = "M2"
$ws.Cells.Item(16,5).Value = [double]"3"
$ws.Cells.Item(16,6).Value = [double]"1"
$ws.Cells.Item(16,7).Value = [double]"21.38374133333333"
$ws.Cells.Item(16,8).Value = [double]"64.151224"
$ws.Cells.Item(16,9).Value = [double]"0.3528756426099016"
$ws.Cells.Item(16,10).Value = [double]"0.3528756426099016"
$ws.Cells.Item(16,11).Value = [double]"3"
$ws.Cells.Item(16,12).Value = [double]"1"
$ws.Cells.Item(16,13).Value = [double]"4.77305"
$ws.Cells.Item(16,14).Value = [double]"14.31915"
$ws.Cells.Item(16,15).Value = [double]"0.03571136528892854"
$ws.Cells.Item(16,16).Value = [double]"0.03571136528892854"
$ws.Cells.Item(16,17).Value = [double]"102.0656665710667"
$ws.Cells.Item(16,18).Value = [double]"918.5909991396002"
$ws.Cells.Item(16,19).Value = [double]"0.01260167097480759"
$ws.Cells.Item(16,20).Value = [double]"0.01260167097480759"

# Row 17
$ws.Cells.Item(17,1).Value = "sCs"
$ws.Cells.Item(17,2).Value = "Sema6d"
$ws.Cells.Item(17,3).Value = "Kdr"
$ws.Cells.Item(17,4).Value = "sCs"
$ws.Cells.Item(17,5).Value = [double]"3"
$ws.Cells.Item(17,6).Value = [double]"1"
$ws.Cells.Item(17,7).Value = [double]"21.38374133333333"
$ws.Cells.Item(17,8).Value = [double]"64.151224"
$ws.Cells.Item(17,9).Value = [double]"0.3528756426099016"
$ws.Cells.Item(17,10).Value = [double]"0.3528756426099016"
$ws.Cells.Item(17,11).Value = [double]"3"
$ws.Cells.Item(17,12).Value = [double]"1"
$ws.Cells.Item(17,13).Value = [double]"0.751186"
$ws.Cells.Item(17,14).Value = [double]"2.253558"
$ws.Cells.Item(17,15).Value = [double]"0.00562028004021099"
$ws.Cells.Item(17,16).Value = [double]"0.00562028004021099"
$ws.Cells.Item(17,17).Value = [double]"16.06316711722133"
$ws.Cells.Item(17,18).Value = [double]"144.568504054992"
$ws.Cells.Item(17,19).Value = [double]"0.001983259930837057"
$ws.Cells.Item(17,20).Value = [double]"0.001983259930837057"
